# Updates cryptos list figures (price / 1h volume change) per the
# upstream GitHub Actions refresh job, plus a couple of rows that swapped
# their rank position (WrappedEther/Polkadot, THORChain->PaxDollar).
#
# Column D ("Price") values that look like plain numbers (e.g. "4.80")
# would otherwise be auto-converted to a numeric type by Excel, silently
# dropping significant trailing zeros / formatting. We force those cells
# to text via a temporary "@" NumberFormat, then ClearFormats() right
# after so the cell ends up with no explicit style, matching the rest of
# the sheet (which never carried per-cell number formats to begin with).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "35.032.58"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +1.06%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.852.59"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +2.34%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "236.94"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +3.23%  "
$ws.Range("E6").Value = "  +1.24%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "42.39"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +7.34%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.329"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +2.79%  "
$ws.Range("E10").Value = "  +2.16%  "
$ws.Range("E11").Value = "  +0.43%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "11.42"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +2.13%  "
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.80"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +5.07%  "
$ws.Range("B15").Value = "WrappedEther"
$ws.Range("C15").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "1.842.81"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +2.11%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.677"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +2.24%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "35.029.92"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +1.12%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "70.27"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +1.50%  "
$ws.Range("E19").Value = "  +2.03%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "240.56"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +0.57%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.19"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +3.23%  "
$ws.Range("E22").Value = "  +3.16%  "
$ws.Range("E23").Value = "  -0.05%  "
$ws.Range("E24").Value = "  +1.06%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "170.55"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -1.59%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.88"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +25.01%  "
$ws.Range("E27").Value = "  +3.81%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.67"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +2.54%  "
$ws.Range("E29").Value = "  +0.51%  "
$ws.Range("E31").Value = "  -0.03%  "
$ws.Range("E32").Value = "  +0.88%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.03"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +3.55%  "
$ws.Range("E34").Value = "  +23.66%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.00"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +12.27%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.33"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +8.34%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.782"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +13.96%  "
$ws.Range("E38").Value = "  +11.68%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0202"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +6.22%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "90.73"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -0.05%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.350.65"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +1.40%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "14.72"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +3.52%  "
$ws.Range("E43").Value = "  +4.16%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "12.78"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +53.61%  "
$ws.Range("E45").Value = "  -0.24%  "
$ws.Range("E46").Value = "  +6.77%  "
$ws.Range("E47").Value = "  +0.22%  "
$ws.Range("E48").Value = "  +6.82%  "
$ws.Range("E49").Value = "  +2.13%  "
$ws.Range("E50").Value = "  +2.72%  "
$ws.Range("B51").Value = "PaxDollar"
$ws.Range("C51").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.01"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -0.04%  "
